$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndicatorCriterion")

# Shift the existing criterion rows 2-11 down to 3-12 (copy cell values only,
# bottom-up so nothing is clobbered before it is read) to make room for the
# new criterion row, matching the dimension growing from A1:E11 to A1:E12.
# Column C ("iri") is left untouched everywhere - it is blank in every data
# row both before and after the edit, so there is nothing to shift there.
foreach ($c in 1, 2, 4, 5) {
    for ($r = 11; $r -ge 2; $r--) {
        $v = $ws.Cells.Item($r, $c).Value()
        $ws.Cells.Item($r + 1, $c).Value = $v
    }
}

# Populate the new criterion row: IntegrationFoodSystemFramework, under the
# CoherenceOnFramework category.
$ws.Cells.Item(2, 1).Value = 'CoherenceOnFramework'
$ws.Cells.Item(2, 2).Value = 'IntegrationFoodSystemFramework'
$ws.Cells.Item(2, 4).Value = 'integration on food system framework'
$ws.Cells.Item(2, 5).Value = 'Integration on food system framework Indicators should be mapped within a food system framework and assessed as part of an integrated set, ensuring each provides a unique, non-redundant contribution. Collectively, they must offer a balanced and comprehensive representation of the food system. The operative allocation involves mapping indicators within the framework and maintaining system-level coherence to guarantee their relevance and effectiveness.'

# The three criteria that used to sit under CoherenceOnFramework (now on rows
# 3-5 after the shift) are recategorised to GoalOrientedAndPolicyRelevant.
$ws.Cells.Item(3, 1).Value = "GoalOrientedAndPolicyRelevant"
$ws.Cells.Item(4, 1).Value = "GoalOrientedAndPolicyRelevant"
$ws.Cells.Item(5, 1).Value = "GoalOrientedAndPolicyRelevant"
